$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: 'Days of Chunder' / 'Antidote'
$ws.Range("H6").Value = 3240.6667
$ws.Range("I6").Value = 3240.6667
$ws.Range("K6").Value = 9722.000100000001
$ws.Range("M6").Value = -9610.000100000001

# Row 129: 'Practical Command' / "Commanding Craftsman's Draught"
$ws.Range("H129").Value = 296321.2
$ws.Range("I129").Value = 1085.4667
$ws.Range("J129").Value = 529402
$ws.Range("K129").Value = 3256.4001
$ws.Range("L129").Value = 1588206
$ws.Range("M129").Value = 1743.5999
$ws.Range("N129").Value = -1598206

# Row 132: 'Fast-forwarding Flora' / 'Growth Formula Lambda'
$ws.Range("H132").Value = 3027709.5
$ws.Range("I132").Value = 628227.75
$ws.Range("J132").Value = 22223562
$ws.Range("K132").Value = 1884683.25
$ws.Range("L132").Value = 66670686
$ws.Range("M132").Value = -1882153.25
$ws.Range("N132").Value = -66675746

# Row 134: 'Binding Spells' / 'Crocodileskin Index'
$ws.Range("H134").Value = 49620
$ws.Range("J134").Value = 49620
$ws.Range("L134").Value = 49620
$ws.Range("N134").Value = -59760

# Row 138: 'All-night Crafting' / "Cunning Craftsman's Tisane"
$ws.Range("H138").Value = 2331.894
$ws.Range("I138").Value = 1898.6342
$ws.Range("J138").Value = 2735.6135
$ws.Range("K138").Value = 5695.902599999999
$ws.Range("L138").Value = 8206.8405
$ws.Range("M138").Value = -555.9025999999994
$ws.Range("N138").Value = -18486.8405

# Row 141: 'Remedy for Reason' / 'Grade 1 Gemdraught of Mind'
$ws.Range("H141").Value = 1761.4286
$ws.Range("I141").Value = 1636.3158
$ws.Range("J141").Value = 2950
$ws.Range("K141").Value = 4908.9474
$ws.Range("L141").Value = 8850
$ws.Range("M141").Value = 271.0526
$ws.Range("N141").Value = -19210

$ws = $wb.Worksheets.Item("ARM")
# Row 6: "Don't Hit Me One More Time" / 'Bronze Hoplon'
$ws.Range("H6").Value = 8000
$ws.Range("I6").Value = 8000
$ws.Range("K6").Value = 8000
$ws.Range("M6").Value = -7827

# Row 13: 'Get into Their Heads' / 'Bronze Chain Coif'
$ws.Range("H13").Value = 1500300
$ws.Range("J13").Value = 600
$ws.Range("L13").Value = 600
$ws.Range("N13").Value = -888

# Row 37: 'Get Shirty' / 'Steel Chainmail'
$ws.Range("H37").Value = 9910.736999999999
$ws.Range("I37").Value = 5333.3335
$ws.Range("K37").Value = 5333.3335
$ws.Range("M37").Value = -5060.3335

$ws = $wb.Worksheets.Item("BSM")
# Row 26: 'Unseamly Conditions' / 'Iron Pickaxe'
$ws.Range("H26").Value = 12100
$ws.Range("I26").Value = 12100
$ws.Range("K26").Value = 12100
$ws.Range("M26").Value = -11808

$ws = $wb.Worksheets.Item("CRP")
# Row 13: 'Compulsory Conjury' / 'Maple Cane'
$ws.Range("H13").Value = 20500
$ws.Range("I13").Value = 20500
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 20500
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -20361
$ws.Range("N13").ClearContents()

# Row 50: 'The Arsenal of Theocracy' / 'Cobalt Halberd'
$ws.Range("H50").Value = 12339.8
$ws.Range("J50").Value = 12339.8
$ws.Range("L50").Value = 12339.8
$ws.Range("N50").Value = -13589.8

# Row 51: 'Greenstone for Greenhorns' / 'Jade Crook'
$ws.Range("H51").Value = 9170.666999999999
$ws.Range("J51").Value = 9170.666999999999
$ws.Range("L51").Value = 9170.666999999999
$ws.Range("N51").Value = -10642.667

# Row 58: 'You Do the Heavy Lifting' / 'Mahogany Lumber'
$ws.Range("H58").Value = 1169092.1
$ws.Range("I58").Value = 4001.9033
$ws.Range("J58").Value = 5683817
$ws.Range("K58").Value = 4001.9033
$ws.Range("L58").Value = 5683817
$ws.Range("M58").Value = -3798.9033
$ws.Range("N58").Value = -5684223

# Row 59: 'Bow Down to Magic' / 'Crab Bow'
$ws.Range("H59").Value = 15796.75
$ws.Range("J59").Value = 15796.75
$ws.Range("L59").Value = 15796.75
$ws.Range("N59").Value = -18086.75

# Row 60: 'Bowing to Greater Power' / 'Yew Longbow'
$ws.Range("H60").Value = 7359.2856
$ws.Range("J60").Value = 8202.5
$ws.Range("L60").Value = 8202.5
$ws.Range("N60").Value = -9224.5

# Row 61: 'Incant Now, Think Later' / 'Jade Crook'
$ws.Range("H61").Value = 9170.666999999999
$ws.Range("J61").Value = 9170.666999999999
$ws.Range("L61").Value = 9170.666999999999
$ws.Range("N61").Value = -9866.666999999999

# Row 68: 'Do You Even String Bow' / 'Holy Cedar Composite Bow'
$ws.Range("H68").Value = 18285.857
$ws.Range("J68").Value = 18285.857
$ws.Range("L68").Value = 18285.857
$ws.Range("N68").Value = -19783.857

# Row 70: 'A Reward Fitting of the Faithful' / 'Holy Cedar Necklace'
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 71: 'Win One Bow, Get Three Free (L)' / 'Holy Cedar Composite Bow'
$ws.Range("H71").Value = 18285.857
$ws.Range("J71").Value = 18285.857
$ws.Range("L71").Value = 54857.571
$ws.Range("N71").Value = -62345.571

# Row 73: 'Just Rewards for Just Devotion (L)' / 'Holy Cedar Necklace'
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 74: 'License to Heal' / 'Dark Chestnut Rod'
$ws.Range("H74").Value = 20138.25
$ws.Range("J74").Value = 20138.25
$ws.Range("L74").Value = 20138.25
$ws.Range("N74").Value = -21886.25

# Row 77: 'Purified Polyrhythm (L)' / 'Dark Chestnut Rod'
$ws.Range("H77").Value = 20138.25
$ws.Range("J77").Value = 20138.25
$ws.Range("L77").Value = 60414.75
$ws.Range("N77").Value = -69150.75

# Row 132: 'Hull Lotta Damage' / 'Ginseng Lumber'
$ws.Range("H132").Value = 1525
$ws.Range("I132").Value = 1055.381
$ws.Range("J132").Value = 2511.2
$ws.Range("K132").Value = 3166.143
$ws.Range("L132").Value = 7533.599999999999
$ws.Range("M132").Value = -636.143
$ws.Range("N132").Value = -12593.6

# Row 136: 'Turali Quality' / 'Dark Mahogany Lumber'
$ws.Range("H136").Value = 1169092.1
$ws.Range("I136").Value = 4001.9033
$ws.Range("J136").Value = 5683817
$ws.Range("K136").Value = 12005.7099
$ws.Range("L136").Value = 17051451
$ws.Range("M136").Value = -9455.7099
$ws.Range("N136").Value = -17056551

$ws = $wb.Worksheets.Item("CUL")
# Row 3: 'Trout Fishing in Limsa' / 'Grilled Trout'
$ws.Range("H3").Value = 8043.3
$ws.Range("I3").Value = 7266.6665
$ws.Range("K3").Value = 21799.9995
$ws.Range("M3").Value = -21687.9995

# Row 11: 'Putting the Squeeze On' / 'Orange Juice'
$ws.Range("H11").Value = 99
$ws.Range("I11").Value = 48.75
$ws.Range("J11").Value = 300
$ws.Range("K11").Value = 146.25
$ws.Range("L11").Value = 900
$ws.Range("M11").Value = -6.25
$ws.Range("N11").Value = -1180

# Row 60: 'Drinking to Your Health' / 'Mulled Tea'
$ws.Range("H60").Value = 6181.857
$ws.Range("I60").Value = 589
$ws.Range("J60").Value = 10376.5
$ws.Range("K60").Value = 1767
$ws.Range("L60").Value = 31129.5
$ws.Range("M60").Value = -1516
$ws.Range("N60").Value = -31631.5

# Row 121: 'A Cookie for Your Troubles' / 'Coffee Biscuit'
$ws.Range("H121").Value = 3125729.2
$ws.Range("I121").Value = 290
$ws.Range("J121").Value = 3255956
$ws.Range("K121").Value = 870
$ws.Range("L121").Value = 9767868
$ws.Range("M121").Value = 440
$ws.Range("N121").Value = -9770488

# Row 131: 'The Mountain Steeped' / 'Tsai tou Vounou'
$ws.Range("H131").Value = 944
$ws.Range("I131").Value = 583.3333
$ws.Range("J131").Value = 955.15466
$ws.Range("K131").Value = 1749.9999
$ws.Range("L131").Value = 2865.46398
$ws.Range("M131").Value = 3290.0001
$ws.Range("N131").Value = -12945.46398

# Row 137: 'Creative Chocolate' / 'Gateau au Chocolat'
$ws.Range("H137").Value = 7344.087
$ws.Range("J137").Value = 10428.143
$ws.Range("L137").Value = 31284.429
$ws.Range("N137").Value = -41484.429

$ws = $wb.Worksheets.Item("GSM")
# Row 70: 'Sky Is the Limit' / 'Mythrite Ingot'
$ws.Range("H70").Value = 5137987.5
$ws.Range("I70").Value = 2319073.5
$ws.Range("J70").Value = 10992656
$ws.Range("K70").Value = 2319073.5
$ws.Range("L70").Value = 10992656
$ws.Range("M70").Value = -2318803.5
$ws.Range("N70").Value = -10993196

# Row 73: 'Hulls of Broken Dreams (L)' / 'Mythrite Ingot'
$ws.Range("H73").Value = 5137987.5
$ws.Range("I73").Value = 2319073.5
$ws.Range("J73").Value = 10992656
$ws.Range("K73").Value = 2319073.5
$ws.Range("L73").Value = 10992656
$ws.Range("M73").Value = -2318137.5
$ws.Range("N73").Value = -10994528

$ws = $wb.Worksheets.Item("LTW")
# Row 16: 'Saddle Sore' / 'Hard Leather'
$ws.Range("H16").Value = 44412.348
$ws.Range("I16").Value = 72115.78999999999
$ws.Range("J16").Value = 1318.1111
$ws.Range("K16").Value = 72115.78999999999
$ws.Range("L16").Value = 1318.1111
$ws.Range("M16").Value = -71945.78999999999
$ws.Range("N16").Value = -1658.1111

# Row 122: 'Hell on Leather' / 'Gaja Leather'
$ws.Range("H122").Value = 16253153
$ws.Range("I122").Value = 11305591
$ws.Range("J122").Value = 40001450
$ws.Range("K122").Value = 33916773
$ws.Range("L122").Value = 120004350
$ws.Range("M122").Value = -33914323
$ws.Range("N122").Value = -120009250

# Row 132: 'Tenets of Tanning' / 'Silver Lobo Leather'
$ws.Range("H132").Value = 5716000.5
$ws.Range("I132").Value = 6212435.5
$ws.Range("J132").Value = 6999.5
$ws.Range("K132").Value = 18637306.5
$ws.Range("L132").Value = 20998.5
$ws.Range("M132").Value = -18634776.5
$ws.Range("N132").Value = -26058.5

# Row 136: "Respect for Br'aax" / "Br'aax Leather"
$ws.Range("H136").Value = 3380366.8
$ws.Range("I136").Value = 5001824
$ws.Range("J136").Value = 2330.4167
$ws.Range("K136").Value = 15005472
$ws.Range("L136").Value = 6991.250100000001
$ws.Range("M136").Value = -15002922
$ws.Range("N136").Value = -12091.2501

$ws = $wb.Worksheets.Item("WVR")
# Row 10: 'Just for Kecks' / 'Hempen Kecks'
$ws.Range("H10").Value = 1000
$ws.Range("I10").Value = 1000
$ws.Range("K10").Value = 1000
$ws.Range("M10").Value = -831

# Row 132: 'Comfy Cabins' / 'Snow Cotton Cloth'
$ws.Range("H132").Value = 868428.4399999999
$ws.Range("I132").Value = 2483.7334
$ws.Range("J132").Value = 1907562.1
$ws.Range("K132").Value = 7451.2002
$ws.Range("L132").Value = 5722686.300000001
$ws.Range("M132").Value = -4921.2002
$ws.Range("N132").Value = -5727746.300000001

# Row 136: 'Weaving the Envelope' / 'Sarcenet Cloth'
$ws.Range("H136").Value = 1979.7188
$ws.Range("I136").Value = 1591.2941
$ws.Range("J136").Value = 2419.9333
$ws.Range("K136").Value = 4773.8823
$ws.Range("L136").Value = 7259.7999
$ws.Range("M136").Value = -2223.8823
$ws.Range("N136").Value = -12359.7999
